$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The BF column holds a "Date" text column that was entered as plain text
# (e.g. "6-27-2012-13"). The fix re-writes it in ISO form ("2013-06-27").
# A bare Value assignment of an ISO-looking string gets auto-coerced into a
# real date serial by the engine, so we force it to stay literal text with a
# leading apostrophe (classic Excel "treat as text" prefix), then reset the
# cell style back to Normal so the quote-prefix formatting doesn't leave a
# stray style on the cell.
for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Cells.Item($row, 58)  # column BF = 58
    $cell.Value = "'2013-06-27"
    $cell.Style = "Normal"
}
